$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the data range to Text format so numeric-looking strings
# (e.g. "1.005", "25.805.35") are preserved as text, matching the
# original inline-string cell types. Reset the style afterwards so
# no stray style index is left on the cells.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '25.805.35'
$ws.Range("E2").Value = '  +6.92%  '
$ws.Range("D3").Value = '1.759.74'
$ws.Range("E3").Value = '  +5.26%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.62%  '
$ws.Range("D5").Value = '316.43'
$ws.Range("E5").Value = '  +2.79%  '
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("D7").Value = '0.3814'
$ws.Range("E7").Value = '  +2.59%  '
$ws.Range("D8").Value = '0.3600'
$ws.Range("E8").Value = '  +4.37%  '
$ws.Range("D9").Value = '50.32'
$ws.Range("E9").Value = '  +4.84%  '
$ws.Range("D10").Value = '1.225'
$ws.Range("E10").Value = '  +4.41%  '
$ws.Range("D11").Value = '0.07681'
$ws.Range("E11").Value = '  +5.80%  '
$ws.Range("D12").Value = '0.9992'
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").Value = '21.64'
$ws.Range("E13").Value = '  +5.37%  '
$ws.Range("E14").Value = '  +7.07%  '
$ws.Range("D15").Value = '7.078'
$ws.Range("E15").Value = '  +5.08%  '
$ws.Range("D16").Value = '1.763.57'
$ws.Range("E16").Value = '  +6.12%  '
$ws.Range("D17").Value = '0.00001153'
$ws.Range("E17").Value = '  +4.70%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = '0.06787'
$ws.Range("E18").Value = '  +1.02%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = '0.9992'
$ws.Range("E19").Value = '  +0.35%  '
$ws.Range("D20").Value = '86.87'
$ws.Range("E20").Value = '  +5.94%  '
$ws.Range("D21").Value = '17.66'
$ws.Range("E21").Value = '  +7.38%  '
$ws.Range("D22").Value = '6.504'
$ws.Range("E22").Value = '  +6.07%  '
$ws.Range("D23").Value = '13.02'
$ws.Range("E23").Value = '  +8.44%  '
$ws.Range("D24").Value = '25.813.62'
$ws.Range("E24").Value = '  +7.32%  '
$ws.Range("D25").Value = '2.441'
$ws.Range("D26").Value = '2.898'
$ws.Range("E26").Value = '  +8.49%  '
$ws.Range("D27").Value = '20.78'
$ws.Range("E27").Value = '  +6.42%  '
$ws.Range("D28").Value = '155.92'
$ws.Range("E28").Value = '  +2.55%  '
$ws.Range("D29").Value = '1.959.45'
$ws.Range("E29").Value = '  +6.26%  '
$ws.Range("D30").Value = '133.82'
$ws.Range("E30").Value = '  +5.31%  '
$ws.Range("D31").Value = '1.207'
$ws.Range("E31").Value = '  +22.95%  '
$ws.Range("D32").Value = '7.191'
$ws.Range("E32").Value = '  +12.84%  '
$ws.Range("D33").Value = '4.211'
$ws.Range("E33").Value = '  +3.39%  '
$ws.Range("D34").Value = '14.28'
$ws.Range("E34").Value = '  +15.91%  '
$ws.Range("D35").Value = '1.807'
$ws.Range("E35").Value = '  +4.43%  '
$ws.Range("D36").Value = '0.08758'
$ws.Range("E36").Value = '  +4.53%  '
$ws.Range("D37").Value = '5.737'
$ws.Range("E37").Value = '  +7.78%  '
$ws.Range("D38").Value = '0.06757'
$ws.Range("E38").Value = '  +6.24%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.02498'
$ws.Range("E39").Value = '  +7.68%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '9.376'
$ws.Range("E40").Value = '  +4.74%  '
$ws.Range("D41").Value = '0.2265'
$ws.Range("E41").Value = '  +8.92%  '
$ws.Range("D42").Value = '1.296'
$ws.Range("E42").Value = '  +0.47%  '
$ws.Range("D43").Value = '0.6591'
$ws.Range("E43").Value = '  +8.04%  '
$ws.Range("D44").Value = '14.40'
$ws.Range("E44").Value = '  +8.51%  '
$ws.Range("D45").Value = '0.9988'
$ws.Range("E45").Value = '  +0.42%  '
$ws.Range("D46").Value = '0.6343'
$ws.Range("E46").Value = '  +6.70%  '
$ws.Range("D47").Value = '3.908'
$ws.Range("E47").Value = '  +2.62%  '
$ws.Range("D48").Value = '2.174'
$ws.Range("E48").Value = '  +8.41%  '
$ws.Range("D49").Value = '132.18'
$ws.Range("E49").Value = '  +3.88%  '
$ws.Range("E50").Value = '  +5.66%  '
$ws.Range("D51").Value = '80.88'
$ws.Range("E51").Value = '  +6.54%  '

# Reset style back to Normal/General so cell styling matches the original.
$dataRange.Style = "Normal"
